$d = $word.ActiveDocument
$d.Content.Find.Execute("MEDISONIC MODELO H60", $true, $false, $false, $false, $false, $true, 1, $false, "MINDRAY MODELO DC – N3", 2)
